{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// This reproduces the diff:\n//  1. Splits the \"Ans: Possible reasons could be:\" paragraph into two runs\n//     (\"Ans\" and \": Possible reasons could be:\"), bracketing the first run\n//     with <w:proofErr w:type=\"spellStart\"/> / <w:proofErr w:type=\"spellEnd\"/>.\n//  2. Expands the JDBC-connection sentence into several additional runs/\n//     sentences describing Hibernate session-level caching, and removes the\n//     trailing \"_GoBack\" bookmark from that paragraph.\n//  3. Adds two new bullet (ListParagraph) items after it (\"Check if any UI\n//     bugs...\" and \"Check if the first entered number...\").\n//  4. Moves the \"_GoBack\" bookmark onto the (still empty) bullet paragraph\n//     that follows those two new items.\n//\n// Because none of this is expressible with plain text/paragraph insertion\n// APIs (proofErr marks, precise run splits, bookmark relocation), the two\n// affected regions are replaced in place using Range.insertOoxml with a\n// \"flat OPC\" package \u2014 the supported way to inject literal WordprocessingML\n// through Office.js.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst WORDML_NS = \"http://schemas.openxmlformats.org/wordprocessingml/2006/main\";\n\nfunction flatOpc(bodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"' + WORDML_NS + '\">' +\n    \"<w:body>\" +\n    bodyXml +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\n// ---------------------------------------------------------------------\n// 1) \"Ans: Possible reasons could be:\" paragraph -> split run + proofErr\n// ---------------------------------------------------------------------\nconst ansParagraph = paragraphs.items.find(\n  (p) => p.text.trim() === \"Ans: Possible reasons could be:\"\n);\nif (!ansParagraph) {\n  throw new Error(\"Could not find the 'Ans: Possible reasons could be:' paragraph\");\n}\n\nconst ansXml =\n  \"<w:p>\" +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>Ans</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  \"<w:r><w:t>: Possible reasons could be:</w:t></w:r>\" +\n  \"</w:p>\";\n\nansParagraph.getRange().insertOoxml(flatOpc(ansXml), \"Replace\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) \"Session might also be the reason...\" paragraph (JDBC sentence) and\n//    the trailing, empty bookmark paragraph right after it -> rewritten\n//    JDBC/Hibernate text (no bookmark) + two new bullet items + the\n//    bookmark relocated onto a new trailing (still empty) bullet item.\n// ---------------------------------------------------------------------\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nconst sessionParagraph = paragraphs2.items.find((p) =>\n  p.text.startsWith(\"Session might also be the reason\")\n);\nif (!sessionParagraph) {\n  throw new Error(\"Could not find the 'Session might also be the reason...' paragraph\");\n}\nconst nextParagraph = sessionParagraph.getNext();\nnextParagraph.load(\"text\");\nawait context.sync();\n\n// Build a Range spanning from the start of the \"Session...\" paragraph to\n// the end of the following (empty, bookmark-holding) paragraph, so both\n// paragraphs get replaced together by the four target paragraphs.\nconst spanRange = sessionParagraph.getRange().expandTo(nextParagraph.getRange());\n\nconst listNumPr =\n  '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>';\n\nconst sessionXml =\n  \"<w:p>\" +\n  listNumPr +\n  \"<w:r><w:t>Session might also be the reason, if the application is following some old architecture.</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> Data might have stored in the session. While saving data to the database, somehow it might have saved the data from the old session.</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">JDBC connection might not have closed, because of why, session might not have flushed during the previous update. </w:t></w:r>' +\n  \"<w:r><w:t>Some</w:t></w:r>\" +\n  \"<w:r><w:t>t</w:t></w:r>\" +\n  \"<w:r><w:t>imes if the object\\u2019s values changed in the DB hiber</w:t></w:r>\" +\n  \"<w:r><w:t>n</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\">ate needs to refresh the session level cache or it might use old values. This is because if the item is in the session, it\\u2019ll pull it from the session instead of the DB. </w:t></w:r>' +\n  \"<w:r><w:t>So it is taking</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> from the previous session only- So for this, in Hibernate</w:t></w:r>' +\n  \"<w:r><w:t>, close the JBDC connection or try a session refresh on your object.</w:t></w:r>\" +\n  \"</w:p>\";\n\nconst uiBugXml =\n  \"<w:p>\" +\n  listNumPr +\n  '<w:r><w:t xml:space=\"preserve\">Check if any UI bugs are there or the form </w:t></w:r>' +\n  \"<w:r><w:t>is working correctly.</w:t></w:r>\" +\n  \"</w:p>\";\n\nconst firstNumberXml =\n  \"<w:p>\" +\n  listNumPr +\n  \"<w:r><w:t>Check if the first entered number is associated with any other account or if there is a requirement saying \\u201cOne number has to be associated with one account only or something like that\\u201d.</w:t></w:r>\" +\n  \"</w:p>\";\n\nconst bookmarkOnlyXml =\n  \"<w:p>\" +\n  listNumPr +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  \"</w:p>\";\n\nspanRange.insertOoxml(\n  flatOpc(sessionXml + uiBugXml + firstNumberXml + bookmarkOnlyXml),\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument / $d is the open document.\n#\n# This reproduces the diff:\n#  1. Splits the \"Ans: Possible reasons could be:\" paragraph into two runs\n#     (\"Ans\" and \": Possible reasons could be:\"), bracketing the first run\n#     with <w:proofErr w:type=\"spellStart\"/> / <w:proofErr w:type=\"spellEnd\"/>.\n#  2. Expands the JDBC-connection sentence into several additional runs/\n#     sentences describing Hibernate session-level caching, and removes the\n#     trailing \"_GoBack\" bookmark from that paragraph.\n#  3. Adds two new bullet (ListParagraph) items after it (\"Check if any UI\n#     bugs...\" and \"Check if the first entered number...\").\n#  4. Moves the \"_GoBack\" bookmark onto the (still empty) bullet paragraph\n#     that follows those two new items.\n#\n# None of this (proofErr marks, precise run splits, bookmark relocation) is\n# expressible with plain Range.Text assignment, so the two affected regions\n# are replaced in place with Range.InsertXML \u2014 the WordprocessingML-literal\n# insertion method on the Range object.\n\n$d = $word.ActiveDocument\n\n# -----------------------------------------------------------------------\n# 1) \"Ans: Possible reasons could be:\" paragraph -> split run + proofErr\n# -----------------------------------------------------------------------\n$ansParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"Ans: Possible reasons could be:\") {\n        $ansParagraph = $p\n        break\n    }\n}\nif ($ansParagraph -eq $null) {\n    throw \"Could not find the 'Ans: Possible reasons could be:' paragraph\"\n}\n\n$ansXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n    '<w:proofErr w:type=\"spellStart\"/>' + `\n    '<w:r><w:t>Ans</w:t></w:r>' + `\n    '<w:proofErr w:type=\"spellEnd\"/>' + `\n    '<w:r><w:t>: Possible reasons could be:</w:t></w:r>' + `\n    '</w:p>'\n\n$ansParagraph.Range.InsertXML($ansXml)\n\n# -----------------------------------------------------------------------\n# 2) \"Session might also be the reason...\" paragraph (JDBC sentence) and\n#    the trailing, empty bookmark paragraph right after it -> rewritten\n#    JDBC/Hibernate text (no bookmark) + two new bullet items + the\n#    bookmark relocated onto a new trailing (still empty) bullet item.\n# -----------------------------------------------------------------------\n$sessionParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"Session might also be the reason\")) {\n        $sessionParagraph = $p\n        break\n    }\n}\nif ($sessionParagraph -eq $null) {\n    throw \"Could not find the 'Session might also be the reason...' paragraph\"\n}\n$nextParagraph = $sessionParagraph.Next()\n\n# Range spanning from the start of the \"Session...\" paragraph to the end of\n# the following (empty, bookmark-holding) paragraph, so both paragraphs get\n# replaced together by the four target paragraphs below.\n$spanRange = $d.Range($sessionParagraph.Range.Start, $nextParagraph.Range.End)\n\n$listNumPr = '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>'\n\n# Only the first paragraph in the concatenated fragment needs the w:\n# namespace declaration; InsertXML accepts a run of sibling <w:p> elements\n# without a synthetic wrapper root.\n$sessionXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + $listNumPr + `\n    '<w:r><w:t>Session might also be the reason, if the application is following some old architecture.</w:t></w:r>' + `\n    '<w:r><w:t xml:space=\"preserve\"> Data might have stored in the session. While saving data to the database, somehow it might have saved the data from the old session.</w:t></w:r>' + `\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' + `\n    '<w:r><w:t xml:space=\"preserve\">JDBC connection might not have closed, because of why, session might not have flushed during the previous update. </w:t></w:r>' + `\n    '<w:r><w:t>Some</w:t></w:r>' + `\n    '<w:r><w:t>t</w:t></w:r>' + `\n    '<w:r><w:t>imes if the object' + [char]0x2019 + 's values changed in the DB hiber</w:t></w:r>' + `\n    '<w:r><w:t>n</w:t></w:r>' + `\n    '<w:r><w:t xml:space=\"preserve\">ate needs to refresh the session level cache or it might use old values. This is because if the item is in the session, it' + [char]0x2019 + 'll pull it from the session instead of the DB. </w:t></w:r>' + `\n    '<w:r><w:t>So it is taking</w:t></w:r>' + `\n    '<w:r><w:t xml:space=\"preserve\"> from the previous session only- So for this, in Hibernate</w:t></w:r>' + `\n    '<w:r><w:t>, close the JBDC connection or try a session refresh on your object.</w:t></w:r>' + `\n    '</w:p>'\n\n$uiBugXml = '<w:p>' + $listNumPr + `\n    '<w:r><w:t xml:space=\"preserve\">Check if any UI bugs are there or the form </w:t></w:r>' + `\n    '<w:r><w:t>is working correctly.</w:t></w:r>' + `\n    '</w:p>'\n\n$firstNumberXml = '<w:p>' + $listNumPr + `\n    '<w:r><w:t>Check if the first entered number is associated with any other account or if there is a requirement saying ' + [char]0x201C + 'One number has to be associated with one account only or something like that' + [char]0x201D + '.</w:t></w:r>' + `\n    '</w:p>'\n\n$bookmarkOnlyXml = '<w:p>' + $listNumPr + `\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' + `\n    '<w:bookmarkEnd w:id=\"0\"/>' + `\n    '</w:p>'\n\n$replacementXml = $sessionXml + $uiBugXml + $firstNumberXml + $bookmarkOnlyXml\n\n$spanRange.InsertXML($replacementXml)\n"}
